$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the website cell (B10) with the new domain
$ws.Range("B10").Value = "www.stat.gov.kg"

# Update the active selection to reflect the edited cell
$ws.Range("B10").Select()
